$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" "60.561.20"
Set-CellText $ws "E2" "  +2.50%  "
Set-CellText $ws "D3" "2.699.87"
Set-CellText $ws "E3" "  +2.40%  "
Set-CellText $ws "D4" "1.00"
Set-CellText $ws "E4" "  +0.04%  "
Set-CellText $ws "D5" "525.11"
Set-CellText $ws "E5" "  +0.35%  "
Set-CellText $ws "D6" "145.32"
Set-CellText $ws "E6" "  -0.45%  "
Set-CellText $ws "D7" "0.996"
Set-CellText $ws "E7" "  +0.02%  "
Set-CellText $ws "D8" "0.576"
Set-CellText $ws "E8" "  +0.42%  "
Set-CellText $ws "D9" "2.729.61"
Set-CellText $ws "E9" "  +2.75%  "
Set-CellText $ws "D10" "6.74"
Set-CellText $ws "E10" "  +6.45%  "
Set-CellText $ws "D11" "0.105"
Set-CellText $ws "E11" "  +0.62%  "
Set-CellText $ws "E12" "  +0.54%  "
Set-CellText $ws "D13" "0.131"
Set-CellText $ws "E13" "  +3.15%  "
Set-CellText $ws "D14" "3.178.70"
Set-CellText $ws "E14" "  +2.51%  "
Set-CellText $ws "D15" "60.572.50"
Set-CellText $ws "E15" "  +2.52%  "
Set-CellText $ws "D16" "21.24"
Set-CellText $ws "E16" "  +1.24%  "
Set-CellText $ws "D17" "2.755.34"
Set-CellText $ws "E17" "  +4.03%  "
Set-CellText $ws "E18" "  +0.39%  "
Set-CellText $ws "D19" "345.10"
Set-CellText $ws "E19" "  -0.51%  "
Set-CellText $ws "E20" "  -0.13%  "
Set-CellText $ws "E21" "  +3.07%  "
Set-CellText $ws "D22" "6.44"
Set-CellText $ws "E22" "  +4.07%  "
Set-CellText $ws "D23" "0.997"
Set-CellText $ws "E23" "  -0.08%  "
Set-CellText $ws "D24" "63.30"
Set-CellText $ws "E24" "  +2.23%  "
Set-CellText $ws "D25" "0.420"
Set-CellText $ws "E25" "  +0.72%  "
Set-CellText $ws "E26" "  +2.40%  "
Set-CellText $ws "E27" "  -0.16%  "
Set-CellText $ws "D28" "0.0₃0819"
Set-CellText $ws "E28" "  +1.92%  "
Set-CellText $ws "D29" "7.26"
Set-CellText $ws "E29" "  +2.01%  "
Set-CellText $ws "D30" "6.82"
Set-CellText $ws "E30" "  +8.71%  "
Set-CellText $ws "D31" "0.998"
Set-CellText $ws "E31" "  -0.02%  "
Set-CellText $ws "E32" "  +0.73%  "
Set-CellText $ws "D33" "19.03"
Set-CellText $ws "E33" "  +0.23%  "
Set-CellText $ws "D34" "149.60"
Set-CellText $ws "E34" "  -0.67%  "
Set-CellText $ws "E35" "  +6.80%  "
Set-CellText $ws "D36" "1.23"
Set-CellText $ws "E36" "  +7.96%  "
Set-CellText $ws "D37" "0.934"
Set-CellText $ws "E37" "  -4.95%  "
Set-CellText $ws "D38" "0.873"
Set-CellText $ws "E38" "  +2.87%  "
Set-CellText $ws "E39" "  +6.87%  "
Set-CellText $ws "D40" "37.11"
Set-CellText $ws "E40" "  +0.85%  "
Set-CellText $ws "E41" "  -0.37%  "
Set-CellText $ws "D42" "281.42"
Set-CellText $ws "E42" "  +0.89%  "
Set-CellText $ws "D43" "20.05"
Set-CellText $ws "E43" "  +2.25%  "
Set-CellText $ws "B44" "Mantle"
Set-CellText $ws "C44" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-CellText $ws "D44" "0.612"
Set-CellText $ws "E44" "  +0.03%  "
Set-CellText $ws "E45" "  +0.02%  "
Set-CellText $ws "B46" "Maker"
Set-CellText $ws "C46" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-CellText $ws "D46" "2.145.03"
Set-CellText $ws "E46" "  +7.39%  "
Set-CellText $ws "D47" "0.0985"
Set-CellText $ws "E47" "  -0.21%  "
Set-CellText $ws "E48" "  +4.82%  "
Set-CellText $ws "E49" "  +2.66%  "
Set-CellText $ws "D50" "10.52"
Set-CellText $ws "E50" "  +2.11%  "
Set-CellText $ws "D51" "0.0231"
Set-CellText $ws "E51" "  +0.71%  "
